$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (Late/heading/Outstanding shift right),
# matching column M's width for the newly inserted column.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make the "Repayment schedule" sheet the active tab and set its selection.
$ws.Activate() | Out-Null
$ws.Range("R7").Select() | Out-Null
